# Update Fgf13-Scn8a NATMI LR-pair data with new TPM-derived values.
# Sending clusters change from {Inflammatory-Mac, MuSCs, Resolving-Mac}
# to {ECs, Inflammatory-Mac, MuSCs, Neutrophils, Resolving-Mac}; target
# clusters change from {ECs, FAPs, MuSCs, Neutrophils, Resolving-Mac}
# to {ECs, FAPs, MuSCs}. All numeric columns (E:T) get refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Fgf13", "Scn8a", "ECs", 1.0, 0.5, 0.1113035, 0.222607, 0.0289261410371415, 0.0273813389619527, 1.0, 0.5, 0.1493335, 0.298667, 0.1277186509504481, 0.1039080998824773, 0.01662134121725, 0.066485364869, 0.003694407710466108, 0.002845142903774548),
    @("ECs", "Fgf13", "Scn8a", "FAPs", 1.0, 0.5, 0.1113035, 0.222607, 0.0289261410371415, 0.0273813389619527, 3.0, 1.0, 0.535862, 1.607586, 0.4583001920909173, 0.5592891302275516, 0.05964331611699999, 0.357859896702, 0.01325685599377092, 0.01531408525249629),
    @("ECs", "Fgf13", "Scn8a", "MuSCs", 1.0, 0.5, 0.1113035, 0.222607, 0.0289261410371415, 0.0273813389619527, 2.0, 1.0, 0.4840425, 0.968085, 0.4139811569586346, 0.3368027698899712, 0.05387562439875, 0.215502497595, 0.01197487733290448, 0.009222110805681857),
    @("Inflammatory-Mac", "Fgf13", "Scn8a", "ECs", 1.0, 0.3333333333333333, 0.406345, 1.219035, 0.1056030832789379, 0.1499450176386367, 1.0, 0.5, 0.1493335, 0.298667, 0.1277186509504481, 0.1039080998824773, 0.0606809210575, 0.364085526345, 0.01348748333259377, 0.01558050186967529),
    @("Inflammatory-Mac", "Fgf13", "Scn8a", "FAPs", 1.0, 0.3333333333333333, 0.406345, 1.219035, 0.1056030832789379, 0.1499450176386367, 3.0, 1.0, 0.535862, 1.607586, 0.4583001920909173, 0.5592891302275516, 0.21774484439, 1.95970359951, 0.04839791335213037, 0.08386261849706803),
    @("Inflammatory-Mac", "Fgf13", "Scn8a", "MuSCs", 1.0, 0.3333333333333333, 0.406345, 1.219035, 0.1056030832789379, 0.1499450176386367, 2.0, 1.0, 0.4840425, 0.968085, 0.4139811569586346, 0.3368027698899712, 0.1966882496625, 1.180129497975, 0.04371768659421375, 0.05050189727189344),
    @("MuSCs", "Fgf13", "Scn8a", "ECs", 2.0, 1.0, 3.3023715, 6.604743, 0.8582377352557337, 0.8124035041107618, 1.0, 0.5, 0.1493335, 0.298667, 0.1277186509504481, 0.1039080998824773, 0.49315469439525, 1.972618777581, 0.1096129657416301, 0.08441530445001559),
    @("MuSCs", "Fgf13", "Scn8a", "FAPs", 2.0, 1.0, 3.3023715, 6.604743, 0.8582377352557337, 0.8124035041107618, 3.0, 1.0, 0.535862, 1.607586, 0.4583001920909173, 0.5592891302275516, 1.769615396733, 10.617692380398, 0.3933305189273765, 0.4543684492079231),
    @("MuSCs", "Fgf13", "Scn8a", "MuSCs", 2.0, 1.0, 3.3023715, 6.604743, 0.8582377352557337, 0.8124035041107618, 2.0, 1.0, 0.4840425, 0.968085, 0.4139811569586346, 0.3368027698899712, 1.59848815678875, 6.393952627155, 0.3552942505867269, 0.2736197504528232),
    @("Neutrophils", "Fgf13", "Scn8a", "ECs", 1.0, 0.3333333333333333, 0.0004946666666666667, 0.001484, 0.0001285565841718604, 0.0001825365196042254, 1.0, 0.5, 0.1493335, 0.298667, 0.1277186509504481, 0.1039080998824773, 0.00007387030466666667, 0.000443221828, 0.00001641907350122773, 0.00001896702291123563),
    @("Neutrophils", "Fgf13", "Scn8a", "FAPs", 1.0, 0.3333333333333333, 0.0004946666666666667, 0.001484, 0.0001285565841718604, 0.0001825365196042254, 3.0, 1.0, 0.535862, 1.607586, 0.4583001920909173, 0.5592891302275516, 0.0002650730693333333, 0.002385657624, 0.0000589175072205158, 0.0001020906912842116),
    @("Neutrophils", "Fgf13", "Scn8a", "MuSCs", 1.0, 0.3333333333333333, 0.0004946666666666667, 0.001484, 0.0001285565841718604, 0.0001825365196042254, 2.0, 1.0, 0.4840425, 0.968085, 0.4139811569586346, 0.3368027698899712, 0.00023943969, 0.00143663814, 0.00005322000345011686, 0.00006147880540877815),
    @("Resolving-Mac", "Fgf13", "Scn8a", "ECs", 1.0, 0.3333333333333333, 0.027337, 0.082011, 0.007104483844015123, 0.01008760276904456, 1.0, 0.5, 0.1493335, 0.298667, 0.1277186509504481, 0.1039080998824773, 0.004082329889500001, 0.024493979337, 0.0009073750922568651, 0.001048183636100637),
    @("Resolving-Mac", "Fgf13", "Scn8a", "FAPs", 1.0, 0.3333333333333333, 0.027337, 0.082011, 0.007104483844015123, 0.01008760276904456, 3.0, 1.0, 0.535862, 1.607586, 0.4583001920909173, 0.5592891302275516, 0.014648859494, 0.131839735446, 0.003255986310418949, 0.005641886578779973),
    @("Resolving-Mac", "Fgf13", "Scn8a", "MuSCs", 1.0, 0.3333333333333333, 0.027337, 0.082011, 0.007104483844015123, 0.01008760276904456, 2.0, 1.0, 0.4840425, 0.968085, 0.4139811569586346, 0.3368027698899712, 0.0132322698225, 0.079393618935, 0.002941122441339308, 0.003397532554163952)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}

Write-Output "updated $($data.Length) rows x $($data[0].Length) cols"
